# The workbook was re-uploaded after a light cleanup pass in Excel:
#   1. The "Good" scenario sheets were renamed to "Biomass".
#   2. A stray, value-less (but style-formatted) trailing row was removed
#      from the bottom of the StockDict_Carbon sheet.
#   3. The first sheet was left as the active tab when the file was saved.

$wb = $excel.ActiveWorkbook

# 1. Rename sheets: FlowDict_Good -> FlowDict_Biomass, StockDict_Good -> StockDict_Biomass
$wb.Worksheets.Item("FlowDict_Good").Name = "FlowDict_Biomass"
$wb.Worksheets.Item("StockDict_Good").Name = "StockDict_Biomass"

# 2. Remove the extra blank (formatted-only) row 23 at the bottom of StockDict_Carbon
$wb.Worksheets.Item("StockDict_Carbon").Rows.Item(23).Delete()

# 3. Leave the first sheet (FlowDict_Biomass) as the active/selected tab
$wb.Worksheets.Item("FlowDict_Biomass").Activate()
